$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 269
$ws.Range("C3").Value = 166845
$ws.Range("C4").Value = 157746
$ws.Range("C7").Value = 5.45
$ws.Range("C8").Value = 65.16
